# Updates cryptos list values (prices & volume%) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.557.72"
$ws.Range("E2").Value = "  +1.09%  "

# Row 3
$ws.Range("D3").Value = "1.853.20"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.56"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4736"
$ws.Range("E7").Value = "  +0.93%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2744"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06322"
$ws.Range("E9").Value = "  -0.89%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.73"
$ws.Range("E10").Value = "  +9.01%  "

# Row 11
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.830.59"
$ws.Range("E11").Value = "  -1.10%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07450"
$ws.Range("E12").Value = "  +0.39%  "

# Row 13
$ws.Range("E13").Value = "  +1.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.51"
$ws.Range("E14").Value = "  -0.68%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6260"
$ws.Range("E15").Value = "  -0.21%  "

# Row 16
$ws.Range("D16").Value = "30.513.14"
$ws.Range("E16").Value = "  +1.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "244.19"
$ws.Range("E17").Value = "  +7.23%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9999"
$ws.Range("E18").Value = "  -0.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.69"
$ws.Range("E19").Value = "  +0.77%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007338"
$ws.Range("E20").Value = "  +0.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.935"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.924"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("E24").Value = "  -0.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.96"
$ws.Range("E25").Value = "  -2.15%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.98"
$ws.Range("E26").Value = "  +1.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.878"
$ws.Range("E27").Value = "  +0.40%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1019"
$ws.Range("E28").Value = "  -1.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.359"
$ws.Range("E29").Value = "  -1.58%  "

# Row 30
$ws.Range("E30").Value = "  -2.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.827"
$ws.Range("E31").Value = "  -1.19%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04838"
$ws.Range("E32").Value = "  -1.12%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.135"
$ws.Range("E33").Value = "  -1.66%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7018"
$ws.Range("E34").Value = "  -1.27%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.695"
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01901"
$ws.Range("E36").Value = "  +2.45%  "

# Row 37
$ws.Range("E37").Value = "  +1.55%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8742"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.992"
$ws.Range("E39").Value = "  +2.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.83"
$ws.Range("E40").Value = "  +1.63%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9994"
$ws.Range("E41").Value = "  +0.14%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.527"
$ws.Range("E42").Value = "  -0.42%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4058"
$ws.Range("E43").Value = "  -0.62%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.183"
$ws.Range("E44").Value = "  +1.83%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.47"
$ws.Range("E45").Value = "  +3.93%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1211"
$ws.Range("E46").Value = "  +1.84%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.54"
$ws.Range("E47").Value = "  +1.34%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.515"
$ws.Range("E48").Value = "  -1.38%  "

# Row 49
$ws.Range("E49").Value = "  -0.41%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.354"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3668"
$ws.Range("E51").Value = "  +0.19%  "

